$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.639.95"
$ws.Range("E2").Value = "  -0.06%  "
$ws.Range("D3").Value = "'1.598.34"
$ws.Range("E3").Value = "  +0.26%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "'211.57"
$ws.Range("E5").Value = "  -0.23%  "
$ws.Range("E6").Value = "  +0.25%  "
$ws.Range("E8").Value = "  +0.20%  "
$ws.Range("E9").Value = "  -0.03%  "
$ws.Range("D10").Value = "'19.55"
$ws.Range("E10").Value = "  -0.51%  "
$ws.Range("D11").Value = "'0.0836"
$ws.Range("E11").Value = "  +0.32%  "
$ws.Range("D12").Value = "'1.823.38"
$ws.Range("E12").Value = "  +0.33%  "
$ws.Range("D13").Value = "'1.604.71"
$ws.Range("E13").Value = "  +0.41%  "
$ws.Range("D14").Value = "'4.03"
$ws.Range("E14").Value = "  -0.23%  "
$ws.Range("D15").Value = "'0.523"
$ws.Range("E15").Value = "  -0.05%  "
$ws.Range("D16").Value = "'64.93"
$ws.Range("E16").Value = "  -0.67%  "
$ws.Range("D17").Value = "'26.641.03"
$ws.Range("D18").Value = "'0.0₃0738"
$ws.Range("E18").Value = "  +1.08%  "
$ws.Range("D19").Value = "'209.21"
$ws.Range("E19").Value = "  +0.08%  "
$ws.Range("E20").Value = "  +0.15%  "
$ws.Range("D21").Value = "'7.07"
$ws.Range("E21").Value = "  +5.33%  "
$ws.Range("D22").Value = "'4.27"
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("D23").Value = "'2.30"
$ws.Range("E23").Value = "  -0.75%  "
$ws.Range("D24").Value = "'8.90"
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("D25").Value = "'145.73"
$ws.Range("E25").Value = "  -0.08%  "
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("D27").Value = "'7.18"
$ws.Range("E27").Value = "  +0.18%  "
$ws.Range("D28").Value = "'0.115"
$ws.Range("E28").Value = "  +0.62%  "
$ws.Range("D29").Value = "'15.27"
$ws.Range("E29").Value = "  -0.44%  "
$ws.Range("D30").Value = "'0.0509"
$ws.Range("E30").Value = "  +0.78%  "
$ws.Range("E31").Value = "  +0.28%  "
$ws.Range("E32").Value = "  -0.14%  "
$ws.Range("E33").Value = "  +1.10%  "
$ws.Range("D34").Value = "'0.629"
$ws.Range("E34").Value = "  -5.97%  "
$ws.Range("D35").Value = "'1.279.81"
$ws.Range("E35").Value = "  -1.76%  "
$ws.Range("E36").Value = "  +0.13%  "
$ws.Range("E37").Value = "  -0.13%  "
$ws.Range("E38").Value = "  -0.43%  "
$ws.Range("D39").Value = "'0.844"
$ws.Range("E39").Value = "  +1.50%  "
$ws.Range("D40").Value = "'5.51"
$ws.Range("E40").Value = "  +2.88%  "
$ws.Range("E41").Value = "  +1.10%  "
$ws.Range("E42").Value = "  -0.33%  "
$ws.Range("D43").Value = "'64.05"
$ws.Range("E43").Value = "  +1.20%  "
$ws.Range("D44").Value = "'0.941"
$ws.Range("E44").Value = "  +17.14%  "
$ws.Range("D45").Value = "'1.736.20"
$ws.Range("E45").Value = "  +0.34%  "
$ws.Range("D46").Value = "'90.29"
$ws.Range("E46").Value = "  +1.24%  "
$ws.Range("E47").Value = "  +0.28%  "
$ws.Range("E48").Value = "  +4.24%  "
$ws.Range("E49").Value = "  +1.21%  "
$ws.Range("D50").Value = "'7.49"
$ws.Range("E50").Value = "  -1.08%  "
$ws.Range("E51").Value = "  +0.21%  "
